$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 header text changes from "input_keyName" to "input_name"
$ws.Range("A1").Value = "input_name"

# New column B header "input_reactSelect3Input", same style as A1 (bold/centered "Pandas" style)
$ws.Range("B1").Value = "input_reactSelect3Input"
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New B2 cell (blank row under header, like A2) so the used range grows to A1:B2
$ws.Range("B2").Style = "Normal"

# Column widths: A -> 12, B -> 25 (ColumnWidth excludes the ~0.8333 char gridline/padding
# that Excel adds to the stored <col width>, so subtract it to land on an exact value)
$ws.Columns.Item(1).ColumnWidth = 11.166666666666666
$ws.Columns.Item(2).ColumnWidth = 24.166666666666668
